$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title of the table (was "Inventory", now "Toys")
$ws.Range("B21").Value = "TableName: Toys, From: Managers&Users via Orders"

# Add new field row for ToyCategory (mirrors existing style of data rows)
$ws.Range("B26").Value = "ToyCategory"
$ws.Range("C26").Value = "Varchar(50)"
$ws.Range("D26").Value = "No"
$ws.Range("E26").Value = "NOT NULL"
$ws.Range("F26").Value = "No"
$ws.Range("G26").Value = "No"
$ws.Range("H26").Value = "NA"

# Add new field row for ToyRating
$ws.Range("B27").Value = "ToyRating"
$ws.Range("C27").Value = "Int(11)"
$ws.Range("D27").Value = "No"

# Add list of toy categories below
$ws.Range("B28").Value = "Age Limit"
$ws.Range("B29").Value = "History"
$ws.Range("B30").Value = "Science"
$ws.Range("B31").Value = "Finance"
$ws.Range("B32").Value = "Logic"
$ws.Range("B33").Value = "Solitaire"
$ws.Range("B34").Value = "Team Play"

# Update the active selection to match the edited state
$ws.Range("C39").Select()
